$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect that it now also covers the 2010-18 baseline.
$ws.Name = "2010 and 2010-18"

# --- Insert a new row 6: a "Baseline_C63_2010-18" result for weather year 2010 ---
$ws.Rows("6").Insert()

$ws.Cells.Item(6,1).Value = "CW3M"
$ws.Cells.Item(6,2).Value = "Baseline_C63_2010-18"
$ws.Cells.Item(6,3).Value = 2010
$ws.Cells.Item(6,4).Value = 1044.2558590000001
$ws.Cells.Item(6,5).Value = 1990.4676509999999
$ws.Cells.Item(6,6).Value = 1.255063
$ws.Cells.Item(6,7).Value = 327.58108499999997
$ws.Cells.Item(6,8).Value = 10.610913999999999
$ws.Cells.Item(6,9).Value = 8.8404570000000007
$ws.Cells.Item(6,10).Value = 814.39868200000001
$ws.Cells.Item(6,11).Value = 93.229797000000005
$ws.Cells.Item(6,12).Value = 1291.7937010000001
$ws.Cells.Item(6,13).Value = 1165.4429929999999
$ws.Cells.Item(6,14).Value = 7166.0473629999997
$ws.Cells.Item(6,15).Value = 29450.638672000001
$ws.Cells.Item(6,16).Value = -0.464943
$ws.Cells.Item(6,17).Value = -0.00013799999999999999
$ws.Cells.Item(6,18).Value = 2010

# row6 needs the same numeric display formats as the neighboring data rows
# (D:M -> 0.00, N:O -> integer, P -> 0.00, Q -> 0.000000).
$ws.Range("D6:M6").NumberFormat = "0.00"
$ws.Range("N6:O6").NumberFormat = "0"
$ws.Range("P6").NumberFormat = "0.00"
$ws.Range("Q6").NumberFormat = "0.000000"

# --- Insert a new (blank, later populated) row 9, below the existing "2010-18" row (now row 8) ---
$ws.Rows("9").Insert()

# Duplicate row 8 ("Baseline_2010_c38 9/6/20", 2010-18) into row 9, then relabel it as
# "Baseline_C63_2010-18" with a marginally updated mass-balance total (O9).
$ws.Range("A8:R8").Copy()
$ws.Range("A9:R9").PasteSpecial()

$ws.Cells.Item(9,2).Value = "Baseline_C63_2010-18"
$ws.Cells.Item(9,15).Value = 27227.338324888889

# Restore the original active-cell selection intent (now pointing at the new row's label cells).
$ws.Range("A6:B6").Select() | Out-Null
